# Update Name of Algo
# Apply numeric corrections to the result_data_RandomForest sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -12.63219999999999
$ws.Range("C4").Value = -14.02719999999999
$ws.Range("E6").Value = 12.02589999999999
$ws.Range("C7").Value = -11.8979
$ws.Range("E7").Value = 12.49739999999999
$ws.Range("C8").Value = -12.291
$ws.Range("E8").Value = 13.28610000000001
$ws.Range("B11").Value = 5.058
$ws.Range("B12").Value = 5.5987
$ws.Range("C12").Value = -14.67700000000002
$ws.Range("C14").Value = -11.643
$ws.Range("B15").Value = 5.275000000000002
$ws.Range("E19").Value = 12.70439999999999
$ws.Range("E21").Value = 12.46009999999998
$ws.Range("C22").Value = -11.38259999999999
$ws.Range("E24").Value = 12.43689999999999
$ws.Range("E25").Value = 13.55840000000001
